$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows that are no longer present (old rows 7-9: Florida Dept.
# of Health / Texas Dept. of Public Safety ID-replacement entries). This shifts
# the old rows 10-14 up to become rows 7-11.
$ws.Rows("7:9").Delete()

# Row 2: WREN DL Replacement -> PAYNE, JOSEPH - DL Replacement
$ws.Range("B2").Value = 81468
$ws.Range("D2").Value = 43074
$ws.Range("F2").Value = "Texas Depar of Public Safety"
$ws.Range("H2").Value = "PAYNE, JOSEPH - DL Replacement"
$ws.Range("J2").Value = ""
$ws.Range("L2").Value = "2242 DL Replacements"
$ws.Range("N2").Value = -11

# Row 3: GOMEZ, ANA MARIA DL Replacement -> DAVIS, JEFFREY BERNARD - Photo ID
$ws.Range("B3").Value = 81469
$ws.Range("D3").Value = 43074
$ws.Range("F3").Value = "Texas Dept. of Public Safety"
$ws.Range("H3").Value = "DAVIS, JEFFREY BERNARD - Photo ID"
$ws.Range("J3").Value = "Ö"
$ws.Range("L3").Value = "2221 Photo IDs"
$ws.Range("N3").Value = -16

# Row 4: THOMPSON, HOWARD LEE DL Replacement -> SANDERS, RHONDA - ID Replacement
$ws.Range("B4").Value = 81470
$ws.Range("D4").Value = 43074
$ws.Range("F4").Value = "Texas Dept. of Public Safety"
$ws.Range("H4").Value = "SANDERS, RHONDA - ID Replacement"
$ws.Range("J4").Value = ""
$ws.Range("L4").Value = "2222 ID Replacements"
$ws.Range("N4").Value = -11

# Row 5: BURTON, JEFFREY Birth Certificate -> FRANCIS, LARRY JOSEPH - Photo ID
$ws.Range("B5").Value = 81471
$ws.Range("D5").Value = 43074
$ws.Range("F5").Value = "Texas Dept. of Public Safety"
$ws.Range("H5").Value = "FRANCIS, LARRY JOSEPH - Photo ID"
$ws.Range("J5").Value = "Ö"
$ws.Range("L5").Value = "2221 Photo IDs"
$ws.Range("N5").Value = -16

# Row 6: BURTON, JEFFREY DL Replacement -> ALLARD, JOHN FRANCIS - Birth Certificate
$ws.Range("B6").Value = 81472
$ws.Range("D6").Value = 43074
$ws.Range("F6").Value = "Vital Statistics Burea City of Houston"
$ws.Range("H6").Value = "ALLARD, JOHN FRANCIS - Birth Certificate"
$ws.Range("J6").Value = "Ö"
$ws.Range("L6").Value = "2261 Houston Birth Certs."
$ws.Range("N6").Value = -23

# Update the selection to match the new workbook state.
$null = $ws.Range("A2:N6").Select()
